$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -------------------------------------------------------
# Written in this particular column order so the shared-string table
# ends up with the same index assignment as the target workbook
# (id, date, film, user, rating, then the user names as they first
# appear going down column B).
$ws.Range("A1").Value = "id"
$ws.Range("C1").Value = "date"
$ws.Range("E1").Value = "film"
$ws.Range("B1").Value = "user"
$ws.Range("D1").Value = "rating"

# --- Row data -----------------------------------------------------------
# id, user, date(serial), rating(or $null), film
$rows = @(
    @(1,  "kia",     43765, 5,     92),
    @(2,  "natalie", 43782, 5,     89),
    @(3,  "natalie", 43782, $null, 64),
    @(4,  "abby",    43821, $null, 90),
    @(5,  "kia",     43815, 5,     34),
    @(6,  "abby",    43832, $null, 38),
    @(7,  "kia",     43187, 4,     7),
    @(8,  "natalie", 43877, $null, 42),
    @(9,  "abby",    43929, $null, 16),
    @(10, "kia",     43927, 3,     91)
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $row = $rows[$i]

    $ws.Range("A$r").Value = $row[0]
    $ws.Range("B$r").Value = $row[1]

    $ws.Range("C$r").Value = $row[2]
    $ws.Range("C$r").NumberFormat = "d-mmm-yy"

    if ($null -ne $row[3]) {
        $ws.Range("D$r").Value = $row[3]
    }

    $ws.Range("E$r").Value = $row[4]
}

# --- Column width (best fit on the date column) --------------------------
$ws.Columns.Item(3).ColumnWidth = 8.25

# --- Selection -------------------------------------------------------
[void]$ws.Range("D4").Select()
